$wb = $excel.ActiveWorkbook

# --- ip_address_list: reorder/update entries, add new ones ---
$ws1 = $wb.Worksheets.Item('ip_address_list')
$ws1.UsedRange.ClearContents()
$ws1.Cells.Item(1,1).Value = '440_Austin'
$ws1.Cells.Item(1,2).Value = '10.96.205.240'
$ws1.Cells.Item(1,3).Value = '255.255.255.0'
$ws1.Cells.Item(1,4).Value = 'FortiClient Austin: 
pass:
1Pm#J@PFIkzM&Q@i 
UVt1@Ex2p78kxp30atD7we@!qGK
FH-2050-20
10.96.205.80'
$ws1.Cells.Item(1,5).Value = 0

$ws1.Cells.Item(2,1).Value = '497_Edcha'
$ws1.Cells.Item(2,2).Value = '172.26.7.240'
$ws1.Cells.Item(2,3).Value = '255.255.255.0'
$ws1.Cells.Item(2,4).Value = 'FortiClient Edcha Ex2p78kxp30'
$ws1.Cells.Item(2,5).Value = 0

$ws1.Cells.Item(3,1).Value = '503_Witte'
$ws1.Cells.Item(3,2).Value = '192.168.0.240'
$ws1.Cells.Item(3,3).Value = '255.255.255.0'
$ws1.Cells.Item(3,4).Value = 'PC:	10.96.205.175
NAS:	10.96.205.166
FH:	10.96.205.154
	10.96.205.267
-----------------------------------------
user:JHV_Vision, omron 
Pass:*Jhv2708
---------------------------------------
FortiClient Austin: 
Pass:
1Pm#J@PFIkzM&Q@i 
UVt1@Ex2p78kxp30atD7we@!qGK'
$ws1.Cells.Item(3,5).Value = 0

$ws1.Cells.Item(4,1).Value = '514_Teleflex'
$ws1.Cells.Item(4,2).Value = '192.168.14.240'
$ws1.Cells.Item(4,3).Value = '255.255.255.0'
$ws1.Cells.Item(4,4).Value = 'PC:192.168.14.240
CAM: 192.168.14.??NAS:192.168.14.245
*******************************
user: Vision
pass: *Jhv2708'
$ws1.Cells.Item(4,5).Value = 1

$ws1.Cells.Item(5,1).Value = '518_Valeo'
$ws1.Cells.Item(5,2).Value = '192.168.208.242'
$ws1.Cells.Item(5,3).Value = '255.255.255.0'
$ws1.Cells.Item(5,5).Value = 0

$ws1.Cells.Item(6,1).Value = '518_Valeo II'
$ws1.Cells.Item(6,2).Value = '192.168.1.243'
$ws1.Cells.Item(6,3).Value = '255.255.255.0'
$ws1.Cells.Item(6,5).Value = 1

$ws1.Cells.Item(7,1).Value = '474 B_Austin'
$ws1.Cells.Item(7,2).Value = '10.96.205.175'
$ws1.Cells.Item(7,3).Value = '255.255.255.0'
$ws1.Cells.Item(7,4).Value = 'PC:	10.96.205.175
NAS:	10.96.205.166
FH:	10.96.205.154
	10.96.205.245
-----------------------------------------
user:JHV_Vision, omron 
Pass:*Jhv2708
---------------------------------------
FortiClient Austin: 
Pass:
1Pm#J@PFIkzM&Q@i 
UVt1@Ex2p78kxp30atD7we@!qGK'
$ws1.Cells.Item(7,5).Value = 0

$ws1.Cells.Item(8,1).Value = 'Domaci Wifi'
$ws1.Cells.Item(8,2).Value = '192.168.1.131'
$ws1.Cells.Item(8,3).Value = '255.255.255.0'
$ws1.Cells.Item(8,5).Value = 0

$ws1.Cells.Item(9,1).Value = '527_Teijin'
$ws1.Cells.Item(9,2).Value = '10.101.28.176'
$ws1.Cells.Item(9,3).Value = '255.255.255.0'
$ws1.Cells.Item(9,4).Value = 'XG-X2900:		10.101.28.175
OP:		10.101.28.117'
$ws1.Cells.Item(9,5).Value = 1

$ws1.Cells.Item(10,1).Value = '515_ZF Stara kkkBoleslav'
$ws1.Cells.Item(10,2).Value = '10.9.250.240'
$ws1.Cells.Item(10,3).Value = '255.255.255.0'
$ws1.Cells.Item(10,4).Value = 'NAS - 10.9.250.100Uer:spravce Pass:Jhv*2708 
User:jhvadmin Pass
123TPV456'
$ws1.Cells.Item(10,5).Value = 1

$ws1.Cells.Item(11,1).Value = '515_ '
$ws1.Cells.Item(11,2).Value = '192.168.000.000'
$ws1.Cells.Item(11,3).Value = '255.255.255.0'
$ws1.Cells.Item(11,5).Value = 1

$ws1.Cells.Item(12,1).Value = '529_Witte'
$ws1.Cells.Item(12,2).Value = '192.168.0.240'
$ws1.Cells.Item(12,3).Value = '255.255.255.0'
$ws1.Cells.Item(12,4).Value = 'Kamera VS-S160MX :192.168.0.18'
$ws1.Cells.Item(12,5).Value = 0

$ws1.Cells.Item(13,1).Value = '511_Teleflex'
$ws1.Cells.Item(13,2).Value = '192.168.1.242'
$ws1.Cells.Item(13,3).Value = '255.255.255.0'
$ws1.Cells.Item(13,4).Value = 'Teleflex '
$ws1.Cells.Item(13,5).Value = $false

# --- ip_address_fav_list: add newly-favourited entries ---
$ws2 = $wb.Worksheets.Item('ip_address_fav_list')
$ws2.Cells.Item(4,1).Value = '515_ZF Stara kkkBoleslav'
$ws2.Cells.Item(4,2).Value = '10.9.250.240'
$ws2.Cells.Item(4,3).Value = '255.255.255.0'
$ws2.Cells.Item(4,4).Value = 'NAS - 10.9.250.100Uer:spravce Pass:Jhv*2708 
User:jhvadmin Pass
123TPV456'
$ws2.Cells.Item(4,5).Value = 1

$ws2.Cells.Item(5,1).Value = '515_ '
$ws2.Cells.Item(5,2).Value = '192.168.000.000'
$ws2.Cells.Item(5,3).Value = '255.255.255.0'
$ws2.Cells.Item(5,5).Value = 1

# --- Settings: reset default interface + startup-disk-view flags ---
$ws4 = $wb.Worksheets.Item('Settings')
$ws4.Cells.Item(1,2).Value = 0
$ws4.Cells.Item(4,2).Value = 0

# --- projects_bin2: move ZF entry to row 1 (un-favourite) and add Teleflex entry ---
$ws5 = $wb.Worksheets.Item('projects_bin2')
$ws5.UsedRange.ClearContents()
$ws5.Cells.Item(1,1).Value = '515_ZF Stara Boleslav'
$ws5.Cells.Item(1,2).Value = '10.9.250.240'
$ws5.Cells.Item(1,3).Value = '255.255.255.0'
$ws5.Cells.Item(1,4).Value = 'NAS - 10.9.250.100Uer:spravce Pass:Jhv*2708 
User:jhvadmin Pass
123TPV456'
$ws5.Cells.Item(1,5).Value = 0

$ws5.Cells.Item(3,1).Value = '511_Teleflex'
$ws5.Cells.Item(3,2).Value = '192.168.1.242'
$ws5.Cells.Item(3,3).Value = '255.255.255.0'
$ws5.Cells.Item(3,4).Value = 'Teleflex '
$ws5.Cells.Item(3,5).Value = $false

$ws5.Cells.Item(4,1).Value = '518_Valeo'
$ws5.Cells.Item(4,2).Value = 'V'
$ws5.Cells.Item(4,3).Value = '\\192.168.208.200\10_vision'
$ws5.Cells.Item(4,4).Value = 'jhv_vision'
$ws5.Cells.Item(4,5).Value = 'Jhv*2708'
$ws5.Cells.Item(4,6).Value = 'první sít, ixon
\\192.168.208.200\10_vision'

Write-Host "edit.ps1 applied"
